# Update the "repaymentstrategy" value on the ProductLoanInput sheet from
# "RBI (India)" to "Overdue/Due Fee/Int,Principal", and reflect the
# resulting view/selection state (Excel naturally moves the active cell /
# scroll position to the cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

$ws.Activate()
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 7
